$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("A7").Value = 9917.93
$ws.Range("B7").Value = 9899.1200000000008
$ws.Range("C7").Value = 78.05
$ws.Range("D7").Value = 78.2
$ws.Range("E7").Value = $false
$ws.Range("F7").Value = 0.19
$ws.Range("G6").Copy($ws.Range("G7"))
$ws.Range("G7").Value = 42613.766331018516
$ws.Range("H7").Value = $true

# Row 8
$ws.Range("A8").Value = 9999.26
$ws.Range("B8").Value = 9917.93
$ws.Range("C8").Value = 77.739999999999995
$ws.Range("D8").Value = 78.38
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = 0.82
$ws.Range("G6").Copy($ws.Range("G8"))
$ws.Range("G8").Value = 42614.674363425926
$ws.Range("H8").Value = $true

# Row 9
$ws.Range("A9").Value = 9945.26
$ws.Range("B9").Value = 9999.26
$ws.Range("C9").Value = 78.36
$ws.Range("D9").Value = 77.94
$ws.Range("E9").Value = $false
$ws.Range("F9").Value = -0.54
$ws.Range("G6").Copy($ws.Range("G9"))
$ws.Range("G9").Value = 42615.752511574072
$ws.Range("H9").Value = $false
